$wb = $excel.ActiveWorkbook

# Hunk 0 (ALC)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H61").Value = 177.75
$ws.Range("I61").Value = 177.75
$ws.Range("K61").Value = 533.25
$ws.Range("M61").Value = -361.25

# Hunk 1 (ALC)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 112855.11
$ws.Range("I70").Value = 0
$ws.Range("J70").Value = 112855.11
$ws.Range("K70").Value = 0
$ws.Range("L70").Value = 338565.33
$ws.Range("M70").ClearContents()
$ws.Range("N70").Value = -339105.33

# Hunk 2 (ALC)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H73").Value = 112855.11
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 112855.11
$ws.Range("K73").Value = 0
$ws.Range("L73").Value = 338565.33
$ws.Range("M73").ClearContents()
$ws.Range("N73").Value = -340437.33

# Hunk 3 (ALC)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H81").Value = 0
$ws.Range("J81").Value = 0
$ws.Range("L81").Value = 0
$ws.Range("N81").ClearContents()

# Hunk 4 (ALC)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H84").Value = 0
$ws.Range("J84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("N84").ClearContents()

# Hunk 5 (ALC)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 2119.6956
$ws.Range("I98").Value = 2029.3334
$ws.Range("J98").Value = 2445
$ws.Range("K98").Value = 2029.3334
$ws.Range("L98").Value = 2445
$ws.Range("M98").Value = -531.3334
$ws.Range("N98").Value = -5441

# Hunk 6 (ALC)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 12499.167
$ws.Range("I116").Value = 3283.5715
$ws.Range("J116").Value = 25401
$ws.Range("K116").Value = 3283.5715
$ws.Range("L116").Value = 25401
$ws.Range("M116").Value = 158.4285
$ws.Range("N116").Value = -32285

# Hunk 7 (ALC)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H122").Value = 2119.6956
$ws.Range("I122").Value = 2029.3334
$ws.Range("J122").Value = 2445
$ws.Range("K122").Value = 6088.0002
$ws.Range("L122").Value = 7335
$ws.Range("M122").Value = -3638.0002
$ws.Range("N122").Value = -12235

# Hunk 8 (ALC)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 1639.7441
$ws.Range("I132").Value = 1667.3572
$ws.Range("J132").Value = 480
$ws.Range("K132").Value = 5002.071599999999
$ws.Range("L132").Value = 1440
$ws.Range("M132").Value = -2472.071599999999
$ws.Range("N132").Value = -6500

# Hunk 9 (ARM)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 901.4
$ws.Range("I97").Value = 1209
$ws.Range("J97").Value = 440
$ws.Range("K97").Value = 1209
$ws.Range("L97").Value = 440
$ws.Range("M97").Value = -713
$ws.Range("N97").Value = -1432

# Hunk 10 (ARM)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 93870.72
$ws.Range("I110").Value = 105501.98
$ws.Range("K110").Value = 105501.98
$ws.Range("M110").Value = -103456.98

# Hunk 11 (BSM)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1135145.6
$ws.Range("I86").Value = 1701827.5
$ws.Range("J86").Value = 1782
$ws.Range("K86").Value = 1701827.5
$ws.Range("L86").Value = 1782
$ws.Range("M86").Value = -1700704.5
$ws.Range("N86").Value = -4028

# Hunk 12 (BSM)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 1135145.6
$ws.Range("I89").Value = 1701827.5
$ws.Range("J89").Value = 1782
$ws.Range("K89").Value = 8509137.5
$ws.Range("L89").Value = 8910
$ws.Range("M89").Value = -8503521.5
$ws.Range("N89").Value = -20142

# Hunk 13 (BSM)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H130").Value = 68750
$ws.Range("J130").Value = 68750
$ws.Range("L130").Value = 68750
$ws.Range("N130").Value = -78790

# Hunk 14 (CRP)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H38").Value = 0
$ws.Range("J38").Value = 0
$ws.Range("L38").Value = 0
$ws.Range("N38").ClearContents()

# Hunk 15 (CRP)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H39").Value = 0
$ws.Range("I39").Value = 0
$ws.Range("K39").Value = 0
$ws.Range("M39").ClearContents()

# Hunk 16 (CRP)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H41").Value = 43040.25
$ws.Range("I41").Value = 20000
$ws.Range("J41").Value = 50720.332
$ws.Range("K41").Value = 20000
$ws.Range("L41").Value = 50720.332
$ws.Range("N41").Value = -51576.332
$ws.Range("M41").Value = -19572

# Hunk 17 (CRP)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H46").Value = 0
$ws.Range("J46").Value = 0
$ws.Range("L46").Value = 0
$ws.Range("N46").ClearContents()

# Hunk 18 (CRP)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H49").Value = 0
$ws.Range("I49").Value = 0
$ws.Range("K49").Value = 0
$ws.Range("M49").ClearContents()

# Hunk 19 (CRP)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H51").Value = 55399.668
$ws.Range("I51").Value = 56900
$ws.Range("J51").Value = 54649.5
$ws.Range("K51").Value = 56900
$ws.Range("L51").Value = 54649.5
$ws.Range("M51").Value = -56164
$ws.Range("N51").Value = -56121.5

# Hunk 20 (CRP)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H55").Value = 0
$ws.Range("I55").Value = 0
$ws.Range("J55").Value = 0
$ws.Range("K55").Value = 0
$ws.Range("L55").Value = 0
$ws.Range("M55").ClearContents()
$ws.Range("N55").ClearContents()

# Hunk 21 (CRP)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H61").Value = 55399.668
$ws.Range("I61").Value = 56900
$ws.Range("J61").Value = 54649.5
$ws.Range("K61").Value = 56900
$ws.Range("L61").Value = 54649.5
$ws.Range("M61").Value = -56552
$ws.Range("N61").Value = -55345.5

# Hunk 22 (CRP)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 4016
$ws.Range("I62").Value = 3372.5
$ws.Range("K62").Value = 3372.5
$ws.Range("M62").Value = -2748.5

# Hunk 23 (CRP)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H65").Value = 4016
$ws.Range("I65").Value = 3372.5
$ws.Range("K65").Value = 16862.5
$ws.Range("M65").Value = -13742.5

# Hunk 24 (CRP)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 325282.84
$ws.Range("I134").Value = 2857.5173
$ws.Range("K134").Value = 8572.5519
$ws.Range("M134").Value = -6037.5519

# Hunk 25 (CUL)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H37").Value = 117765.836
$ws.Range("J37").Value = 117765.836
$ws.Range("L37").Value = 353297.508
$ws.Range("N37").Value = -353521.508

# Hunk 26 (CUL)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H134").Value = 1263.6666
$ws.Range("I134").Value = 1263.6666
$ws.Range("K134").Value = 3790.9998
$ws.Range("M134").Value = 1279.0002

# Hunk 27 (GSM)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H17").Value = 3625.75
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 3625.75
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 3625.75
$ws.Range("M17").ClearContents()
$ws.Range("N17").Value = -3961.75

# Hunk 28 (GSM)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 3614.4783
$ws.Range("I126").Value = 3416.3333
$ws.Range("J126").Value = 3741.8572
$ws.Range("K126").Value = 10248.9999
$ws.Range("L126").Value = 11225.5716
$ws.Range("M126").Value = -7778.999899999999
$ws.Range("N126").Value = -16165.5716

# Hunk 29 (GSM)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 35169.188
$ws.Range("I132").Value = 4454.7915
$ws.Range("K132").Value = 13364.3745
$ws.Range("M132").Value = -10834.3745

# Hunk 30 (LTW)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2174.6
$ws.Range("J22").Value = 4499.6665
$ws.Range("L22").Value = 4499.6665
$ws.Range("N22").Value = -5089.6665

# Hunk 31 (LTW)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 2174.6
$ws.Range("J27").Value = 4499.6665
$ws.Range("L27").Value = 4499.6665
$ws.Range("N27").Value = -4713.6665

# Hunk 32 (LTW)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 1077.5264
$ws.Range("I55").Value = 255.41667
$ws.Range("J55").Value = 2486.8572
$ws.Range("K55").Value = 255.41667
$ws.Range("L55").Value = 2486.8572
$ws.Range("M55").Value = -82.41667000000001
$ws.Range("N55").Value = -2832.8572

# Hunk 33 (LTW)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 5679.7085
$ws.Range("I132").Value = 4733.533
$ws.Range("J132").Value = 7256.6665
$ws.Range("K132").Value = 14200.599
$ws.Range("L132").Value = 21769.9995
$ws.Range("M132").Value = -11670.599
$ws.Range("N132").Value = -26829.9995

# Hunk 34 (LTW)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 913525.4
$ws.Range("I136").Value = 956645.5
$ws.Range("J136").Value = 8003
$ws.Range("K136").Value = 2869936.5
$ws.Range("L136").Value = 24009
$ws.Range("M136").Value = -2867386.5
$ws.Range("N136").Value = -29109

# Hunk 35 (WVR)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 34824.844
$ws.Range("I132").Value = 2428.4546
$ws.Range("K132").Value = 7285.3638
$ws.Range("M132").Value = -4755.3638
